# Auto-generated edit script: reorder Python set literal string representations in column D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D3').Value = '{''chika'', ''$'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range('D4').Value = '{''chika'', ''$'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range('D5').Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range('D13').Value = '{'')'', '','', ''('', ''='', '';''}'
$ws.Range('D14').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D15').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D16').Value = '{''<='', ''//='', ''*='', ''||'', ''to'', ''>='', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''='', ''&&'', ''+'', ''+='', ''step'', ''//'', ''%'', '','', '':'', ''-='', ''**='', ''!='', ''/='', '')'', ''<'', ''%='', ''>'', ''}'', '';''}'
$ws.Range('D17').Value = '{'','', '']'', ''}'', '';''}'
$ws.Range('D18').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D19').Value = '{''<='', ''||'', ''to'', ''>='', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''&&'', ''+'', ''step'', ''//'', ''%'', '','', '':'', ''!='', '')'', ''<'', ''>'', ''}'', '';''}'
$ws.Range('D20').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D21').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D23').Value = '{'','', '';''}'
$ws.Range('D25').Value = '{'','', '';''}'
$ws.Range('D26').Value = '{'','', '';''}'
$ws.Range('D32').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D33').Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range('D37').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D39').Value = '{''step'', '','', '':'', ''to'', '')'', '']'', ''}'', '';''}'
$ws.Range('D40').Value = '{''step'', '','', '':'', ''to'', '')'', '']'', ''}'', '';''}'
$ws.Range('D41').Value = '{''step'', '','', '':'', ''to'', '')'', '']'', ''}'', '';''}'
$ws.Range('D42').Value = '{''step'', '','', '':'', ''to'', '')'', '']'', ''}'', '';''}'
$ws.Range('D43').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D44').Value = '{''eme'', ''len'', ''('', ''++'', ''--'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range('D45').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D47').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D48').Value = '{''<='', ''||'', ''to'', ''>='', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''&&'', ''+'', ''step'', ''//'', ''%'', '','', '':'', ''!='', '')'', ''<'', ''>'', ''}'', ''id'', '';''}'
$ws.Range('D49').Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range('D50').Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range('D51').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D52').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D54').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D55').Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''{'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range('D56').Value = '{'')'', '';''}'
$ws.Range('D57').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D60').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D61').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D62').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D64').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D67').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D69').Value = '{''push'', ''ganern'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D70').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D71').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D72').Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''ditech'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range('D75').Value = '{'')'', ''step''}'
$ws.Range('D78').Value = '{'')'', ''step'', ''to''}'
$ws.Range('D80').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D81').Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''ditech'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range('D83').Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range('D84').Value = '{''}'', ''ditech''}'
$ws.Range('D86').Value = '{''amaccana'', ''betsung'', ''ditech'', ''}''}'
$ws.Range('D87').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D88').Value = '{''}'', ''ditech''}'
$ws.Range('D89').Value = '{''ditech'', ''betsung'', ''}''}'
$ws.Range('D90').Value = '{''ditech'', ''betsung'', ''}''}'
$ws.Range('D92').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D93').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range('D94').Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
